$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, with formatting copied from the other
# header cells (G1), matching the bold/centered/bordered header style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add values for the new "Save" column data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
